$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.790.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.21%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.854.02'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.49%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.015'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -2.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.67'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.20%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.014'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.85%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4311'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3790'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.71%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8849'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.13%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.69'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.36%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.860.58'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.50%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.769'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.38%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.490'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.15%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07137'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.13%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.019'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.91%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009036'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.71%  '

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.83%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.52'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.02%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.867.29'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.41%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.282'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.47%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.19'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.78%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.092.90'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.029'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.74%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.40'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.15%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.13%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.045'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.38%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.424'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.22%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.79'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.83%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08968'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.33%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.240'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.19%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7803'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.06%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.581'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.03%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.930'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.29%  '

# Row 36
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.015'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.97%  '

# Row 37
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.145'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.62%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05319'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.56%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01969'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.91%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.875'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.17%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5201'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.17%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.005'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.19%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.43%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.796'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.04%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '110.63'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.98%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.79'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.52%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4746'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.00%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.712'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.51%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06525'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.49%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.015'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.01%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.882'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.58%  '

